$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2471.9092
$ws.Range("J17").Value = 2517.7673
$ws.Range("L17").Value = 7553.3019
$ws.Range("N17").Value = -7889.3019
$ws.Range("H74").Value = 2992.12
$ws.Range("I74").Value = 3098.5
$ws.Range("J74").Value = 2856.7273
$ws.Range("K74").Value = 3098.5
$ws.Range("L74").Value = 2856.7273
$ws.Range("M74").Value = -2162.5
$ws.Range("N74").Value = -4728.7273
$ws.Range("H77").Value = 2992.12
$ws.Range("I77").Value = 3098.5
$ws.Range("J77").Value = 2856.7273
$ws.Range("K77").Value = 15492.5
$ws.Range("L77").Value = 14283.6365
$ws.Range("M77").Value = -10812.5
$ws.Range("N77").Value = -23643.6365
$ws.Range("H137").Value = 31251934
$ws.Range("I137").Value = 1403.0588
$ws.Range("K137").Value = 4209.1764
$ws.Range("M137").Value = -1659.1764

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H45").Value = 1517.4
$ws.Range("I45").Value = 2248
$ws.Range("J45").Value = 1334.75
$ws.Range("K45").Value = 2248
$ws.Range("L45").Value = 1334.75
$ws.Range("M45").Value = -1871
$ws.Range("N45").Value = -2088.75
$ws.Range("H74").Value = 16134617
$ws.Range("I74").Value = 27778746
$ws.Range("K74").Value = 27778746
$ws.Range("M74").Value = -27777872
$ws.Range("H77").Value = 16134617
$ws.Range("I77").Value = 27778746
$ws.Range("K77").Value = 138893730
$ws.Range("M77").Value = -138889362
$ws.Range("H88").Value = 2331.889
$ws.Range("I88").Value = 2149.5
$ws.Range("K88").Value = 2149.5
$ws.Range("M88").Value = -1743.5
$ws.Range("H91").Value = 2331.889
$ws.Range("I91").Value = 2149.5
$ws.Range("K91").Value = 2149.5
$ws.Range("M91").Value = -745.5
$ws.Range("H132").Value = 912416.4399999999
$ws.Range("I132").Value = 1279748
$ws.Range("J132").Value = 85920.5
$ws.Range("K132").Value = 3839244
$ws.Range("L132").Value = 257761.5
$ws.Range("M132").Value = -3836714
$ws.Range("N132").Value = -262821.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7449421.5
$ws.Range("I134").Value = 10056118
$ws.Range("K134").Value = 30168354
$ws.Range("M134").Value = -30165819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2166.6667
$ws.Range("J10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("N10").Value = -3278
$ws.Range("H19").Value = 1526
$ws.Range("I19").Value = 701.3333
$ws.Range("J19").Value = 4000
$ws.Range("K19").Value = 701.3333
$ws.Range("L19").Value = 4000
$ws.Range("M19").Value = -531.3333
$ws.Range("N19").Value = -4340
$ws.Range("H24").Value = 1526
$ws.Range("I24").Value = 701.3333
$ws.Range("J24").Value = 4000
$ws.Range("K24").Value = 701.3333
$ws.Range("L24").Value = 4000
$ws.Range("M24").Value = -531.3333
$ws.Range("N24").Value = -4340
$ws.Range("H60").Value = 8900.571
$ws.Range("J60").Value = 10860.8
$ws.Range("L60").Value = 10860.8
$ws.Range("N60").Value = -11882.8
$ws.Range("H62").Value = 2740
$ws.Range("I62").Value = 2700
$ws.Range("J62").Value = 2766.6667
$ws.Range("K62").Value = 2700
$ws.Range("L62").Value = 2766.6667
$ws.Range("M62").Value = -2076
$ws.Range("N62").Value = -4014.6667
$ws.Range("H65").Value = 2740
$ws.Range("I65").Value = 2700
$ws.Range("J65").Value = 2766.6667
$ws.Range("K65").Value = 13500
$ws.Range("L65").Value = 13833.3335
$ws.Range("M65").Value = -10380
$ws.Range("N65").Value = -20073.3335
$ws.Range("H74").Value = 13588.8
$ws.Range("I74").Value = 3142.5
$ws.Range("K74").Value = 3142.5
$ws.Range("M74").Value = -2268.5
$ws.Range("H77").Value = 13588.8
$ws.Range("I77").Value = 3142.5
$ws.Range("K77").Value = 9427.5
$ws.Range("M77").Value = -5059.5
$ws.Range("H132").Value = 3210.087
$ws.Range("I132").Value = 3293.1428
$ws.Range("J132").Value = 3080.889
$ws.Range("K132").Value = 9879.428400000001
$ws.Range("L132").Value = 9242.667000000001
$ws.Range("M132").Value = -7349.428400000001
$ws.Range("N132").Value = -14302.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 2597.5
$ws.Range("J54").Value = 2597.5
$ws.Range("L54").Value = 7792.5
$ws.Range("N54").Value = -8910.5
$ws.Range("H131").Value = 719.59
$ws.Range("J131").Value = 764.2024
$ws.Range("L131").Value = 2292.6072
$ws.Range("N131").Value = -12372.6072
$ws.Range("H137").Value = 33194.945
$ws.Range("I137").Value = 3905
$ws.Range("J137").Value = 36856.188
$ws.Range("K137").Value = 11715
$ws.Range("L137").Value = 110568.564
$ws.Range("M137").Value = -6615
$ws.Range("N137").Value = -120768.564

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 2200
$ws.Range("J10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("N10").Value = -3338
$ws.Range("H20").Value = 2800
$ws.Range("I20").Value = 3000
$ws.Range("K20").Value = 3000
$ws.Range("M20").Value = -2755

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J20").Value = 4000
$ws.Range("L20").Value = 4000
$ws.Range("N20").Value = -4452
$ws.Range("H22").Value = 433.33334
$ws.Range("I22").Value = 333.33334
$ws.Range("J22").Value = 633.3333
$ws.Range("K22").Value = 333.33334
$ws.Range("L22").Value = 633.3333
$ws.Range("M22").Value = -38.33334000000002
$ws.Range("N22").Value = -1223.3333
$ws.Range("H23").Value = 3500
$ws.Range("H27").Value = 433.33334
$ws.Range("I27").Value = 333.33334
$ws.Range("J27").Value = 633.3333
$ws.Range("K27").Value = 333.33334
$ws.Range("L27").Value = 633.3333
$ws.Range("M27").Value = -226.33334
$ws.Range("N27").Value = -847.3333
$ws.Range("H82").Value = 1559.08
$ws.Range("J82").Value = 2015.3636
$ws.Range("L82").Value = 2015.3636
$ws.Range("N82").Value = -2737.3636
$ws.Range("H85").Value = 1559.08
$ws.Range("J85").Value = 2015.3636
$ws.Range("L85").Value = 2015.3636
$ws.Range("N85").Value = -4511.3636

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 23000
$ws.Range("J64").Value = 23000
$ws.Range("L64").Value = 23000
$ws.Range("N64").Value = -23496
$ws.Range("H67").Value = 23000
$ws.Range("J67").Value = 23000
$ws.Range("L67").Value = 23000
$ws.Range("N67").Value = -24716
$ws.Range("H81").Value = 3732.8572
$ws.Range("J81").Value = 4922.857
$ws.Range("L81").Value = 9845.714
$ws.Range("N81").Value = -11967.714
$ws.Range("H84").Value = 3732.8572
$ws.Range("J84").Value = 4922.857
$ws.Range("L84").Value = 49228.57
$ws.Range("N84").Value = -59836.57
